# DOMA-3372 - add a `role` field (column G) to the contacts export template.
#
# The template has a header row (row 1, i18n placeholders), a "row i" sample
# row (row 2) and a "row i+1" sample row (row 3) that the real export
# generator duplicates down the sheet. We add a new "role" placeholder after
# "email" in each of those three rows, and drop the old filler rows 4-10
# that only existed to carry extra (unused) border styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new "role" column --------------------------------------
# Clone column F's formatting (font/fill/border) into the new column G for
# the three live rows, and give G the same column width as F.
$ws.Range("F1:F3").Copy() | Out-Null
$ws.Range("G1:G3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("G1").ColumnWidth = $ws.Range("F1").ColumnWidth

# --- 2. Fill in the new cell text ---------------------------------------
$ws.Range("G1").Value = "{d.i18n.role}"
$ws.Range("G2").Value = "{d.contacts[i].role}"
$ws.Range("G3").Value = "{d.contacts[i+1].role}"

# --- 3. Remove the now unused filler rows 4-10 --------------------------
$ws.Rows("4:10").Delete()
